# Rename the "SwateTemplateMetadata" sheet to "isa_template"
$wb = $excel.ActiveWorkbook
$metaSheet = $wb.Worksheets.Item("SwateTemplateMetadata")
$metaSheet.Name = "isa_template"

# Make the renamed sheet the active tab (tabSelected moves from Sheet1 to isa_template)
$metaSheet.Activate()

# Update the selection on the isa_template sheet from E19 to D19
$metaSheet.Range("D19").Select()
